# test ksvmeans voi k = 3
# Add a new "K=3, Aroon=5, No volume" results table to the K-SVMeans sheet,
# add a 6th column (F) of results to the existing K=2 table, and make the
# K-SVMeans sheet the active one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K-SVMeans")

# --- New column F for the existing K=2 table (rows 4-9) ---
$ws.Range("F4").Value = 59.42
$ws.Range("F5").Value = 55.67
$ws.Range("F6").Value = 51.87
$ws.Range("F7").Value = 57.61
$ws.Range("F8").Value = 55.84
$ws.Range("F9").Formula = "=AVERAGE(F4:F8)"

# F9 picks up the same "total" style as the other total cells (B9/C9/E9).
$ws.Range("E9").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null

# F6 matches the highlighted style used on E8 in the same row block.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null

# --- New "K = 3, Aroon = 5, No volume" table (rows 11-18) ---
$ws.Range("A11").Value = "K= 3, Aroon = 5, No volume"

$ws.Range("B12").Value = "period = 1"
$ws.Range("C12").Value = "period = 5"

$ws.Range("A13").Value = "BT6"
$ws.Range("B13").Value = 64.52

$ws.Range("A14").Value = "DHG"
$ws.Range("B14").Value = 59.14

$ws.Range("A15").Value = "FPT"
$ws.Range("B15").Value = 65.95

$ws.Range("A16").Value = "VIS"
$ws.Range("B16").Value = 51.89

$ws.Range("A17").Value = "VNM"
$ws.Range("B17").Value = 56.89

$ws.Range("A18").Value = "Total"
$ws.Range("B18").Formula = "=AVERAGE(B13:B17)"

# B14/B15 pick up the same highlighted style as E8.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("B14:B15").PasteSpecial(-4122) | Out-Null

# B18 picks up the same "total" style as the other total cells.
$ws.Range("E9").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null

# Page setup tweak for the sheet.
$ws.PageSetup.Orientation = 1

# Move the active selection/tab onto K-SVMeans (it becomes the active sheet).
$ws.Range("D15").Select() | Out-Null
$ws.Activate()
